# "Update countries & provincias Spain"
# COVID-19 tracker refresh: bump the "last updated" timestamp, update a
# handful of per-country case counters, and re-sync the small block of
# low-count countries (whose rank shuffled slightly between refreshes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- refresh timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 03:16"

# --- Noruega (row 17) ---
$ws.Range("B17").Value = 2625
$ws.Range("C17").Value = 240
$ws.Range("E17").Value = 2609

# --- Suecia (row 20) ---
$ws.Range("E20").Value = 2003
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 27

# --- Hungria (row 74) ---
$ws.Range("D74").Value = 21
$ws.Range("E74").Value = 138
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 8

# --- Consejo Danes para los Refugiados (row 114) ---
$ws.Range("E114").Value = 34
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 2

# --- low-count country block (rows 152-161) re-synced to the latest ranking ---
$ws.Range("A152").Value = "Benin"

$ws.Range("A153").Value = "Haiti"
$ws.Range("B153").Value = 6
$ws.Range("C153").Value = 4
$ws.Range("E153").Value = 6

$ws.Range("A154").Value = "Gabon"
$ws.Range("B154").Value = 6
$ws.Range("C154").Value = 1
$ws.Range("H154").Value = 1

$ws.Range("A155").Value = "Surinam"
$ws.Range("C155").Value = 0
$ws.Range("E155").Value = 5
$ws.Range("H155").Value = 0

$ws.Range("A156").Value = "Islas Caimanes"
$ws.Range("C156").Value = 2

$ws.Range("A157").Value = "Namibia"

$ws.Range("A160").Value = "Suazilandia"
$ws.Range("C160").Value = 0

$ws.Range("A161").Value = "Congo"
$ws.Range("C161").Value = 1
